$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Range("C12").Value = 4
$ws.Range("D12").Value = 108
$ws.Range("G12").Value = 108
$ws.Range("H12").Value = 2.4

# Row 13
$ws.Range("C13").Value = 10
$ws.Range("D13").Value = 113
$ws.Range("G13").Value = 116
$ws.Range("H13").Value = 2.58

# Row 16
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 48
$ws.Range("G16").Value = 50
$ws.Range("H16").Value = 1.25

# Row 17
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 52
$ws.Range("G17").Value = 57
$ws.Range("H17").Value = 1.42

# Row 18
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 66
$ws.Range("G18").Value = 66
$ws.Range("H18").Value = 1.65
